# Clean and document code
# - Fix a stray shared-string value in F3 (was "0.0006812 ***", now a literal "\")
# - Remove the "lfuel" robustness block (old rows 15-18: California/Michigan/
#   SouthCarolina/Wisconsin), shifting the "avg_temp" block (old rows 19-21:
#   Iowa/SouthCarolina/Wisconsin) up to become rows 15-17
# - Fill in the previously-blank "tci" (G) values for that block
# - Leave the selection where the author last clicked

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the stray estimate-column entry on row 3
$ws.Range("F3").Value = "\"

# Drop the four "lfuel" robustness-check rows entirely, shifting rows below up
$ws.Range("A15:J18").Delete(-4162)

# The "avg_temp" block (now rows 15-17) was missing its tci values; fill them in
$ws.Range("G15").Value = 2
$ws.Range("G16").Value = 2
$ws.Range("G17").Value = 7

# Match the author's final selection
$ws.Range("G18").Select()
